$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "team record" columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$hdr = $ws.Range("AD1:AF1")

# Match formatting used by the other header cells (bold font, thin border,
# centered horizontal/top vertical alignment)
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Fill the team record values for every player row (2-54)
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
